$wb = $excel.ActiveWorkbook

$sheet2 = $wb.Worksheets.Item("Sheet2")
$deviceSetup = $wb.Worksheets.Item("DeviceSetupLogins")
$sheet2.Move($deviceSetup)

$sheet3 = $wb.Worksheets.Item("Sheet3")
$deviceSetup2 = $wb.Worksheets.Item("DeviceSetupLogins")
$sheet3.Move($deviceSetup2)

$wb.Worksheets.Item("DeviceSetupLogins").Name = "DeviceSetupLogins_REMOVE"
$wb.Worksheets.Item("InjectSpecificUser").Name = "InjectSpecificUser_REMOVE"

$wb.Worksheets.Item("Sheet1").Range("B9").Select()
$wb.Worksheets.Item("InjectSpecificUser_REMOVE").Range("A28").Select()
$wb.Worksheets.Item("Sheet1").Activate()

foreach ($ws in $wb.Worksheets) {
    Write-Host ("Sheet: " + $ws.Name)
}
